$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Montenegro Prva Liga")

$ws.Cells.Item(3, 2).Value = 6951073
$ws.Cells.Item(3, 5).Value = 'FK Decic Tuzi'
$ws.Cells.Item(3, 6).Value = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(3, 7).Value = 2
$ws.Cells.Item(3, 8).Value = 1
$ws.Cells.Item(3, 9).Value = 'H'
$ws.Cells.Item(3, 10).Value = 1.727
$ws.Cells.Item(3, 11).Value = 3.5
$ws.Cells.Item(3, 12).Value = 4
$ws.Cells.Item(3, 13).Value = 1.6
$ws.Cells.Item(3, 14).Value = 3.6
$ws.Cells.Item(3, 15).Value = 4.75
$ws.Cells.Item(3, 16).Value = -0.75
$ws.Cells.Item(3, 17).Value = 1.8
$ws.Cells.Item(3, 18).Value = 2
$ws.Cells.Item(3, 20).Value = 1.8
$ws.Cells.Item(3, 21).Value = 2
$ws.Cells.Item(3, 22).Value = 0.6000000000000001
$ws.Cells.Item(3, 23).Value = -1
$ws.Cells.Item(3, 25).Value = 0.4
$ws.Cells.Item(3, 26).Value = -0.5
$ws.Cells.Item(3, 27).Value = 0.8
$ws.Cells.Item(3, 28).Value = -1
$ws.Cells.Item(4, 2).Value = 6951072
$ws.Cells.Item(4, 5).Value = 'OFK Petrovac'
$ws.Cells.Item(4, 6).Value = 'FK Mornar Bar'
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 9).Value = 'D'
$ws.Cells.Item(4, 10).Value = 2.375
$ws.Cells.Item(4, 11).Value = 2.875
$ws.Cells.Item(4, 12).Value = 2.875
$ws.Cells.Item(4, 13).Value = 2.15
$ws.Cells.Item(4, 14).Value = 3
$ws.Cells.Item(4, 15).Value = 3.2
$ws.Cells.Item(4, 16).Value = -0.25
$ws.Cells.Item(4, 17).Value = 1.9
$ws.Cells.Item(4, 18).Value = 1.9
$ws.Cells.Item(4, 20).Value = 1.95
$ws.Cells.Item(4, 21).Value = 1.85
$ws.Cells.Item(4, 22).Value = -1
$ws.Cells.Item(4, 23).Value = 2
$ws.Cells.Item(4, 25).Value = -0.5
$ws.Cells.Item(4, 26).Value = 0.45
$ws.Cells.Item(4, 27).Value = -1
$ws.Cells.Item(4, 28).Value = 0.8500000000000001
$ws.Cells.Item(19, 2).Value = 6815422
$ws.Cells.Item(19, 5).Value = 'OFK Mladost DG'
$ws.Cells.Item(19, 6).Value = 'FK Decic Tuzi'
$ws.Cells.Item(19, 8).Value = 2
$ws.Cells.Item(19, 9).Value = 'A'
$ws.Cells.Item(19, 10).Value = 2.4
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 2.75
$ws.Cells.Item(19, 13).Value = 3.1
$ws.Cells.Item(19, 14).Value = 3
$ws.Cells.Item(19, 15).Value = 2.15
$ws.Cells.Item(19, 16).Value = 0.25
$ws.Cells.Item(19, 17).Value = 1.875
$ws.Cells.Item(19, 18).Value = 1.925
$ws.Cells.Item(19, 20).Value = 2.025
$ws.Cells.Item(19, 21).Value = 1.775
$ws.Cells.Item(19, 23).Value = -1
$ws.Cells.Item(19, 24).Value = 1.15
$ws.Cells.Item(19, 25).Value = -1
$ws.Cells.Item(19, 26).Value = 0.925
$ws.Cells.Item(19, 27).Value = -0.5
$ws.Cells.Item(19, 28).Value = 0.3875
$ws.Cells.Item(21, 2).Value = 6815304
$ws.Cells.Item(21, 5).Value = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(21, 6).Value = 'Sutjeska Niksic'
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 'D'
$ws.Cells.Item(21, 10).Value = 5.5
$ws.Cells.Item(21, 11).Value = 3.75
$ws.Cells.Item(21, 12).Value = 1.5
$ws.Cells.Item(21, 13).Value = 3.6
$ws.Cells.Item(21, 14).Value = 3.2
$ws.Cells.Item(21, 15).Value = 1.909
$ws.Cells.Item(21, 16).Value = 0.5
$ws.Cells.Item(21, 17).Value = 1.825
$ws.Cells.Item(21, 18).Value = 1.975
$ws.Cells.Item(21, 20).Value = 1.875
$ws.Cells.Item(21, 21).Value = 1.925
$ws.Cells.Item(21, 23).Value = 2.2
$ws.Cells.Item(21, 24).Value = -1
$ws.Cells.Item(21, 25).Value = 0.825
$ws.Cells.Item(21, 26).Value = -1
$ws.Cells.Item(21, 27).Value = -1
$ws.Cells.Item(21, 28).Value = 0.925
$ws.Cells.Item(29, 2).Value = 6815311
$ws.Cells.Item(29, 5).Value = 'FK Mornar Bar'
$ws.Cells.Item(29, 6).Value = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(29, 7).Value = 4
$ws.Cells.Item(29, 8).Value = 3
$ws.Cells.Item(29, 9).Value = 'H'
$ws.Cells.Item(29, 10).Value = 1.833
$ws.Cells.Item(29, 11).Value = 3.1
$ws.Cells.Item(29, 12).Value = 4
$ws.Cells.Item(29, 13).Value = 2.25
$ws.Cells.Item(29, 14).Value = 2.9
$ws.Cells.Item(29, 15).Value = 3.1
$ws.Cells.Item(29, 16).Value = -0.25
$ws.Cells.Item(29, 17).Value = 1.975
$ws.Cells.Item(29, 18).Value = 1.825
$ws.Cells.Item(29, 19).Value = 2
$ws.Cells.Item(29, 20).Value = 1.9
$ws.Cells.Item(29, 21).Value = 1.9
$ws.Cells.Item(29, 22).Value = 1.25
$ws.Cells.Item(29, 24).Value = -1
$ws.Cells.Item(29, 25).Value = 0.9750000000000001
$ws.Cells.Item(29, 26).Value = -1
$ws.Cells.Item(29, 27).Value = 0.8999999999999999
$ws.Cells.Item(29, 28).Value = -1
$ws.Cells.Item(30, 2).Value = 6815312
$ws.Cells.Item(30, 5).Value = 'Buducnost Podgorica'
$ws.Cells.Item(30, 6).Value = 'FK Arsenal'
$ws.Cells.Item(30, 7).Value = 2
$ws.Cells.Item(30, 8).Value = 0
$ws.Cells.Item(30, 10).Value = 1.444
$ws.Cells.Item(30, 11).Value = 4
$ws.Cells.Item(30, 12).Value = 6
$ws.Cells.Item(30, 13).Value = 1.4
$ws.Cells.Item(30, 14).Value = 4
$ws.Cells.Item(30, 15).Value = 6.5
$ws.Cells.Item(30, 16).Value = -1.25
$ws.Cells.Item(30, 17).Value = 1.95
$ws.Cells.Item(30, 18).Value = 1.85
$ws.Cells.Item(30, 19).Value = 2.5
$ws.Cells.Item(30, 20).Value = 1.775
$ws.Cells.Item(30, 21).Value = 1.925
$ws.Cells.Item(30, 22).Value = 0.3999999999999999
$ws.Cells.Item(30, 25).Value = 0.95
$ws.Cells.Item(30, 27).Value = -1
$ws.Cells.Item(30, 28).Value = 0.925
$ws.Cells.Item(31, 2).Value = 6815315
$ws.Cells.Item(31, 5).Value = 'FK Decic Tuzi'
$ws.Cells.Item(31, 6).Value = 'FK Rudar Pljevlja'
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 1
$ws.Cells.Item(31, 9).Value = 'A'
$ws.Cells.Item(31, 10).Value = 1.615
$ws.Cells.Item(31, 11).Value = 3.5
$ws.Cells.Item(31, 12).Value = 4.75
$ws.Cells.Item(31, 14).Value = 3.8
$ws.Cells.Item(31, 17).Value = 2
$ws.Cells.Item(31, 18).Value = 1.8
$ws.Cells.Item(31, 20).Value = 1.95
$ws.Cells.Item(31, 21).Value = 1.85
$ws.Cells.Item(31, 22).Value = -1
$ws.Cells.Item(31, 24).Value = 5.5
$ws.Cells.Item(31, 25).Value = -1
$ws.Cells.Item(31, 26).Value = 0.8
$ws.Cells.Item(31, 28).Value = 0.8500000000000001
$ws.Cells.Item(54, 2).Value = 6815334
$ws.Cells.Item(54, 5).Value = 'Sutjeska Niksic'
$ws.Cells.Item(54, 6).Value = 'FK Mornar Bar'
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 1
$ws.Cells.Item(54, 9).Value = 'A'
$ws.Cells.Item(54, 10).Value = 1.444
$ws.Cells.Item(54, 11).Value = 4
$ws.Cells.Item(54, 12).Value = 6.5
$ws.Cells.Item(54, 13).Value = 1.444
$ws.Cells.Item(54, 14).Value = 4
$ws.Cells.Item(54, 15).Value = 6.5
$ws.Cells.Item(54, 16).Value = -1.25
$ws.Cells.Item(54, 17).Value = 2
$ws.Cells.Item(54, 18).Value = 1.8
$ws.Cells.Item(54, 19).Value = 2.5
$ws.Cells.Item(54, 20).Value = 2
$ws.Cells.Item(54, 21).Value = 1.8
$ws.Cells.Item(54, 22).Value = -1
$ws.Cells.Item(54, 24).Value = 5.5
$ws.Cells.Item(54, 25).Value = -1
$ws.Cells.Item(54, 26).Value = 0.8
$ws.Cells.Item(54, 27).Value = -1
$ws.Cells.Item(54, 28).Value = 0.8
$ws.Cells.Item(55, 2).Value = 6815426
$ws.Cells.Item(55, 5).Value = 'FK Decic Tuzi'
$ws.Cells.Item(55, 6).Value = 'Buducnost Podgorica'
$ws.Cells.Item(55, 7).Value = 2
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 9).Value = 'H'
$ws.Cells.Item(55, 10).Value = 3
$ws.Cells.Item(55, 11).Value = 3
$ws.Cells.Item(55, 12).Value = 2.25
$ws.Cells.Item(55, 13).Value = 2.75
$ws.Cells.Item(55, 14).Value = 3
$ws.Cells.Item(55, 15).Value = 2.4
$ws.Cells.Item(55, 16).Value = 0.25
$ws.Cells.Item(55, 17).Value = 1.7
$ws.Cells.Item(55, 18).Value = 2.1
$ws.Cells.Item(55, 19).Value = 2.25
$ws.Cells.Item(55, 20).Value = 1.95
$ws.Cells.Item(55, 21).Value = 1.85
$ws.Cells.Item(55, 22).Value = 1.75
$ws.Cells.Item(55, 24).Value = -1
$ws.Cells.Item(55, 25).Value = 0.7
$ws.Cells.Item(55, 26).Value = -1
$ws.Cells.Item(55, 27).Value = -0.5
$ws.Cells.Item(55, 28).Value = 0.425
$ws.Cells.Item(62, 2).Value = 7366684
$ws.Cells.Item(62, 5).Value = 'FK Rudar Pljevlja'
$ws.Cells.Item(62, 6).Value = 'OFK Petrovac'
$ws.Cells.Item(62, 7).Value = 1
$ws.Cells.Item(62, 10).Value = 2.875
$ws.Cells.Item(62, 11).Value = 2.9
$ws.Cells.Item(62, 12).Value = 2.375
$ws.Cells.Item(62, 13).Value = 2.625
$ws.Cells.Item(62, 14).Value = 2.9
$ws.Cells.Item(62, 15).Value = 2.55
$ws.Cells.Item(62, 16).Value = 0
$ws.Cells.Item(62, 17).Value = 1.925
$ws.Cells.Item(62, 18).Value = 1.875
$ws.Cells.Item(62, 19).Value = 2.25
$ws.Cells.Item(62, 20).Value = 1.925
$ws.Cells.Item(62, 21).Value = 1.875
$ws.Cells.Item(62, 22).Value = 1.625
$ws.Cells.Item(62, 25).Value = 0.925
$ws.Cells.Item(62, 28).Value = 0.875
$ws.Cells.Item(64, 2).Value = 6815343
$ws.Cells.Item(64, 5).Value = 'Sutjeska Niksic'
$ws.Cells.Item(64, 6).Value = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(64, 7).Value = 2
$ws.Cells.Item(64, 10).Value = 1.333
$ws.Cells.Item(64, 11).Value = 4.2
$ws.Cells.Item(64, 12).Value = 8
$ws.Cells.Item(64, 13).Value = 1.333
$ws.Cells.Item(64, 14).Value = 4.2
$ws.Cells.Item(64, 15).Value = 8
$ws.Cells.Item(64, 16).Value = -1.5
$ws.Cells.Item(64, 17).Value = 1.975
$ws.Cells.Item(64, 18).Value = 1.825
$ws.Cells.Item(64, 19).Value = 2.75
$ws.Cells.Item(64, 20).Value = 1.9
$ws.Cells.Item(64, 21).Value = 1.9
$ws.Cells.Item(64, 22).Value = 0.333
$ws.Cells.Item(64, 25).Value = 0.9750000000000001
$ws.Cells.Item(64, 28).Value = 0.8999999999999999
$ws.Cells.Item(75, 2).Value = 6815358
$ws.Cells.Item(75, 5).Value = 'OFK Petrovac'
$ws.Cells.Item(75, 6).Value = 'FK Arsenal'
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = 1
$ws.Cells.Item(75, 9).Value = 'D'
$ws.Cells.Item(75, 10).Value = 2.1
$ws.Cells.Item(75, 11).Value = 3.1
$ws.Cells.Item(75, 12).Value = 3.2
$ws.Cells.Item(75, 13).Value = 1.75
$ws.Cells.Item(75, 14).Value = 3.3
$ws.Cells.Item(75, 15).Value = 4.2
$ws.Cells.Item(75, 16).Value = -0.5
$ws.Cells.Item(75, 17).Value = 1.8
$ws.Cells.Item(75, 18).Value = 2
$ws.Cells.Item(75, 19).Value = 2.25
$ws.Cells.Item(75, 20).Value = 1.95
$ws.Cells.Item(75, 21).Value = 1.85
$ws.Cells.Item(75, 22).Value = -1
$ws.Cells.Item(75, 23).Value = 2.3
$ws.Cells.Item(75, 26).Value = 1
$ws.Cells.Item(75, 27).Value = -0.5
$ws.Cells.Item(75, 28).Value = 0.425
$ws.Cells.Item(76, 2).Value = 6815359
$ws.Cells.Item(76, 5).Value = 'Buducnost Podgorica'
$ws.Cells.Item(76, 6).Value = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(76, 7).Value = 3
$ws.Cells.Item(76, 8).Value = 2
$ws.Cells.Item(76, 9).Value = 'H'
$ws.Cells.Item(76, 10).Value = 1.333
$ws.Cells.Item(76, 11).Value = 4.333
$ws.Cells.Item(76, 12).Value = 7.5
$ws.Cells.Item(76, 13).Value = 1.333
$ws.Cells.Item(76, 14).Value = 4.333
$ws.Cells.Item(76, 15).Value = 8
$ws.Cells.Item(76, 16).Value = -1.5
$ws.Cells.Item(76, 17).Value = 1.875
$ws.Cells.Item(76, 18).Value = 1.925
$ws.Cells.Item(76, 19).Value = 2.75
$ws.Cells.Item(76, 20).Value = 1.8
$ws.Cells.Item(76, 21).Value = 2
$ws.Cells.Item(76, 22).Value = 0.333
$ws.Cells.Item(76, 23).Value = -1
$ws.Cells.Item(76, 25).Value = -1
$ws.Cells.Item(76, 26).Value = 0.925
$ws.Cells.Item(76, 27).Value = 0.8
$ws.Cells.Item(76, 28).Value = -1
$ws.Cells.Item(77, 2).Value = 6815357
$ws.Cells.Item(77, 5).Value = 'OFK Mladost DG'
$ws.Cells.Item(77, 6).Value = 'Sutjeska Niksic'
$ws.Cells.Item(77, 10).Value = 4.8
$ws.Cells.Item(77, 11).Value = 3.5
$ws.Cells.Item(77, 12).Value = 1.615
$ws.Cells.Item(77, 13).Value = 4.75
$ws.Cells.Item(77, 14).Value = 3.6
$ws.Cells.Item(77, 15).Value = 1.571
$ws.Cells.Item(77, 16).Value = 1
$ws.Cells.Item(77, 17).Value = 1.775
$ws.Cells.Item(77, 18).Value = 2.025
$ws.Cells.Item(77, 20).Value = 1.825
$ws.Cells.Item(77, 21).Value = 1.975
$ws.Cells.Item(77, 23).Value = 2.6
$ws.Cells.Item(77, 25).Value = 0.7749999999999999
$ws.Cells.Item(77, 26).Value = -1
$ws.Cells.Item(77, 28).Value = 0.4875
$ws.Cells.Item(81, 2).Value = 6815362
$ws.Cells.Item(81, 5).Value = 'Sutjeska Niksic'
$ws.Cells.Item(81, 6).Value = 'FK Decic Tuzi'
$ws.Cells.Item(81, 7).Value = 1
$ws.Cells.Item(81, 8).Value = 1
$ws.Cells.Item(81, 9).Value = 'D'
$ws.Cells.Item(81, 10).Value = 2.2
$ws.Cells.Item(81, 11).Value = 3
$ws.Cells.Item(81, 12).Value = 3.1
$ws.Cells.Item(81, 13).Value = 2.375
$ws.Cells.Item(81, 14).Value = 2.875
$ws.Cells.Item(81, 15).Value = 3
$ws.Cells.Item(81, 16).Value = -0.25
$ws.Cells.Item(81, 17).Value = 2.05
$ws.Cells.Item(81, 18).Value = 1.75
$ws.Cells.Item(81, 19).Value = 2
$ws.Cells.Item(81, 20).Value = 1.8
$ws.Cells.Item(81, 21).Value = 2
$ws.Cells.Item(81, 22).Value = -1
$ws.Cells.Item(81, 23).Value = 1.875
$ws.Cells.Item(81, 26).Value = 0.375
$ws.Cells.Item(81, 27).Value = 0
$ws.Cells.Item(81, 28).Value = 0
$ws.Cells.Item(82, 2).Value = 6815430
$ws.Cells.Item(82, 5).Value = 'Buducnost Podgorica'
$ws.Cells.Item(82, 6).Value = 'FK Mornar Bar'
$ws.Cells.Item(82, 7).Value = 4
$ws.Cells.Item(82, 8).Value = 3
$ws.Cells.Item(82, 9).Value = 'H'
$ws.Cells.Item(82, 10).Value = 1.444
$ws.Cells.Item(82, 11).Value = 3.75
$ws.Cells.Item(82, 12).Value = 6.5
$ws.Cells.Item(82, 13).Value = 1.4
$ws.Cells.Item(82, 14).Value = 4
$ws.Cells.Item(82, 15).Value = 7
$ws.Cells.Item(82, 16).Value = -1.25
$ws.Cells.Item(82, 17).Value = 1.875
$ws.Cells.Item(82, 18).Value = 1.925
$ws.Cells.Item(82, 19).Value = 2.5
$ws.Cells.Item(82, 20).Value = 1.775
$ws.Cells.Item(82, 21).Value = 1.925
$ws.Cells.Item(82, 22).Value = 0.3999999999999999
$ws.Cells.Item(82, 23).Value = -1
$ws.Cells.Item(82, 26).Value = 0.4625
$ws.Cells.Item(82, 27).Value = 0.7749999999999999
$ws.Cells.Item(82, 28).Value = -1
$ws.Cells.Item(105, 2).Value = 6815382
$ws.Cells.Item(105, 5).Value = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(105, 6).Value = 'Sutjeska Niksic'
$ws.Cells.Item(105, 8).Value = 1
$ws.Cells.Item(105, 10).Value = 3.75
$ws.Cells.Item(105, 11).Value = 3.3
$ws.Cells.Item(105, 12).Value = 1.833
$ws.Cells.Item(105, 13).Value = 4.75
$ws.Cells.Item(105, 14).Value = 3
$ws.Cells.Item(105, 15).Value = 1.75
$ws.Cells.Item(105, 16).Value = 0.5
$ws.Cells.Item(105, 17).Value = 1.975
$ws.Cells.Item(105, 18).Value = 1.825
$ws.Cells.Item(105, 19).Value = 2
$ws.Cells.Item(105, 20).Value = 1.85
$ws.Cells.Item(105, 21).Value = 1.95
$ws.Cells.Item(105, 24).Value = 0.75
$ws.Cells.Item(105, 26).Value = 0.825
$ws.Cells.Item(105, 27).Value = -1
$ws.Cells.Item(105, 28).Value = 0.95
$ws.Cells.Item(106, 2).Value = 6815434
$ws.Cells.Item(106, 5).Value = 'OFK Mladost DG'
$ws.Cells.Item(106, 6).Value = 'FK Decic Tuzi'
$ws.Cells.Item(106, 8).Value = 3
$ws.Cells.Item(106, 10).Value = 4.6
$ws.Cells.Item(106, 11).Value = 3.6
$ws.Cells.Item(106, 12).Value = 1.615
$ws.Cells.Item(106, 13).Value = 8
$ws.Cells.Item(106, 14).Value = 4.75
$ws.Cells.Item(106, 15).Value = 1.25
$ws.Cells.Item(106, 16).Value = 0.75
$ws.Cells.Item(106, 17).Value = 1.925
$ws.Cells.Item(106, 18).Value = 1.875
$ws.Cells.Item(106, 19).Value = 2.25
$ws.Cells.Item(106, 20).Value = 1.9
$ws.Cells.Item(106, 21).Value = 1.9
$ws.Cells.Item(106, 24).Value = 0.25
$ws.Cells.Item(106, 26).Value = 0.875
$ws.Cells.Item(106, 27).Value = 0.8999999999999999
$ws.Cells.Item(106, 28).Value = -1
$ws.Cells.Item(107, 2).Value = 7890506
$ws.Cells.Item(107, 5).Value = 'FK Mornar Bar'
$ws.Cells.Item(107, 6).Value = 'FK Arsenal'
$ws.Cells.Item(107, 7).Value = 0
$ws.Cells.Item(107, 8).Value = 0
$ws.Cells.Item(107, 10).Value = 1.85
$ws.Cells.Item(107, 12).Value = 3.9
$ws.Cells.Item(107, 13).Value = 1.85
$ws.Cells.Item(107, 14).Value = 3.3
$ws.Cells.Item(107, 15).Value = 3.5
$ws.Cells.Item(107, 17).Value = 1.925
$ws.Cells.Item(107, 18).Value = 1.875
$ws.Cells.Item(107, 19).Value = 2
$ws.Cells.Item(107, 23).Value = 2.3
$ws.Cells.Item(107, 26).Value = 0.875
$ws.Cells.Item(107, 27).Value = -1
$ws.Cells.Item(107, 28).Value = 0.8500000000000001
$ws.Cells.Item(108, 2).Value = 7890508
$ws.Cells.Item(108, 5).Value = 'OFK Petrovac'
$ws.Cells.Item(108, 6).Value = 'FK Rudar Pljevlja'
$ws.Cells.Item(108, 7).Value = 1
$ws.Cells.Item(108, 8).Value = 1
$ws.Cells.Item(108, 10).Value = 1.75
$ws.Cells.Item(108, 12).Value = 4.5
$ws.Cells.Item(108, 13).Value = 1.8
$ws.Cells.Item(108, 14).Value = 3.2
$ws.Cells.Item(108, 15).Value = 4
$ws.Cells.Item(108, 17).Value = 1.875
$ws.Cells.Item(108, 18).Value = 1.925
$ws.Cells.Item(108, 19).Value = 2.25
$ws.Cells.Item(108, 23).Value = 2.2
$ws.Cells.Item(108, 26).Value = 0.925
$ws.Cells.Item(108, 27).Value = -0.5
$ws.Cells.Item(108, 28).Value = 0.425
$ws.Cells.Item(121, 2).Value = 6815397
$ws.Cells.Item(121, 5).Value = 'FK Arsenal'
$ws.Cells.Item(121, 6).Value = 'OFK Petrovac'
$ws.Cells.Item(121, 8).Value = 1
$ws.Cells.Item(121, 9).Value = 'D'
$ws.Cells.Item(121, 10).Value = 2.6
$ws.Cells.Item(121, 11).Value = 2.7
$ws.Cells.Item(121, 12).Value = 2.8
$ws.Cells.Item(121, 13).Value = 2.75
$ws.Cells.Item(121, 14).Value = 2.5
$ws.Cells.Item(121, 15).Value = 2.875
$ws.Cells.Item(121, 16).Value = 0
$ws.Cells.Item(121, 17).Value = 1.85
$ws.Cells.Item(121, 18).Value = 1.95
$ws.Cells.Item(121, 19).Value = 2
$ws.Cells.Item(121, 20).Value = 2.025
$ws.Cells.Item(121, 21).Value = 1.775
$ws.Cells.Item(121, 23).Value = 1.5
$ws.Cells.Item(121, 24).Value = -1
$ws.Cells.Item(121, 25).Value = 0
$ws.Cells.Item(121, 26).Value = 0
$ws.Cells.Item(121, 27).Value = 0
$ws.Cells.Item(121, 28).Value = 0
$ws.Cells.Item(122, 2).Value = 6815398
$ws.Cells.Item(122, 5).Value = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(122, 6).Value = 'Buducnost Podgorica'
$ws.Cells.Item(122, 8).Value = 2
$ws.Cells.Item(122, 9).Value = 'A'
$ws.Cells.Item(122, 10).Value = 5.75
$ws.Cells.Item(122, 11).Value = 4
$ws.Cells.Item(122, 12).Value = 1.444
$ws.Cells.Item(122, 13).Value = 6.5
$ws.Cells.Item(122, 14).Value = 4.2
$ws.Cells.Item(122, 15).Value = 1.4
$ws.Cells.Item(122, 16).Value = 1.25
$ws.Cells.Item(122, 17).Value = 1.9
$ws.Cells.Item(122, 18).Value = 1.9
$ws.Cells.Item(122, 19).Value = 2.75
$ws.Cells.Item(122, 20).Value = 1.875
$ws.Cells.Item(122, 21).Value = 1.925
$ws.Cells.Item(122, 23).Value = -1
$ws.Cells.Item(122, 24).Value = 0.3999999999999999
$ws.Cells.Item(122, 25).Value = 0.45
$ws.Cells.Item(122, 26).Value = -0.5
$ws.Cells.Item(122, 27).Value = 0.4375
$ws.Cells.Item(122, 28).Value = -0.5
$ws.Cells.Item(126, 2).Value = 6815401
$ws.Cells.Item(126, 5).Value = 'FK Decic Tuzi'
$ws.Cells.Item(126, 6).Value = 'Sutjeska Niksic'
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 9).Value = 'D'
$ws.Cells.Item(126, 10).Value = 2.55
$ws.Cells.Item(126, 12).Value = 2.6
$ws.Cells.Item(126, 13).Value = 2.1
$ws.Cells.Item(126, 14).Value = 3.1
$ws.Cells.Item(126, 15).Value = 3.3
$ws.Cells.Item(126, 16).Value = -0.25
$ws.Cells.Item(126, 17).Value = 1.825
$ws.Cells.Item(126, 18).Value = 1.975
$ws.Cells.Item(126, 19).Value = 2
$ws.Cells.Item(126, 20).Value = 1.925
$ws.Cells.Item(126, 21).Value = 1.875
$ws.Cells.Item(126, 23).Value = 2.1
$ws.Cells.Item(126, 24).Value = -1
$ws.Cells.Item(126, 25).Value = -0.5
$ws.Cells.Item(126, 26).Value = 0.4875
$ws.Cells.Item(126, 28).Value = 0.875
$ws.Cells.Item(127, 2).Value = 6815402
$ws.Cells.Item(127, 5).Value = 'FK Rudar Pljevlja'
$ws.Cells.Item(127, 6).Value = 'FK Jezero'
$ws.Cells.Item(127, 8).Value = 1
$ws.Cells.Item(127, 9).Value = 'A'
$ws.Cells.Item(127, 10).Value = 2.8
$ws.Cells.Item(127, 12).Value = 2.375
$ws.Cells.Item(127, 13).Value = 2.45
$ws.Cells.Item(127, 14).Value = 2.9
$ws.Cells.Item(127, 15).Value = 2.75
$ws.Cells.Item(127, 16).Value = 0
$ws.Cells.Item(127, 17).Value = 1.775
$ws.Cells.Item(127, 18).Value = 2.025
$ws.Cells.Item(127, 19).Value = 1.75
$ws.Cells.Item(127, 20).Value = 1.825
$ws.Cells.Item(127, 21).Value = 1.975
$ws.Cells.Item(127, 23).Value = -1
$ws.Cells.Item(127, 24).Value = 1.75
$ws.Cells.Item(127, 25).Value = -1
$ws.Cells.Item(127, 26).Value = 1.025
$ws.Cells.Item(127, 28).Value = 0.9750000000000001
$ws.Cells.Item(130, 2).Value = 6815405
$ws.Cells.Item(130, 5).Value = 'FK Arsenal'
$ws.Cells.Item(130, 6).Value = 'FK Decic Tuzi'
$ws.Cells.Item(130, 8).Value = 2
$ws.Cells.Item(130, 9).Value = 'A'
$ws.Cells.Item(130, 10).Value = 4
$ws.Cells.Item(130, 11).Value = 3.2
$ws.Cells.Item(130, 12).Value = 1.8
$ws.Cells.Item(130, 13).Value = 5
$ws.Cells.Item(130, 14).Value = 3.25
$ws.Cells.Item(130, 15).Value = 1.65
$ws.Cells.Item(130, 16).Value = 0.75
$ws.Cells.Item(130, 17).Value = 1.875
$ws.Cells.Item(130, 18).Value = 1.925
$ws.Cells.Item(130, 19).Value = 2
$ws.Cells.Item(130, 20).Value = 1.875
$ws.Cells.Item(130, 21).Value = 1.925
$ws.Cells.Item(130, 23).Value = -1
$ws.Cells.Item(130, 24).Value = 0.6499999999999999
$ws.Cells.Item(130, 25).Value = -1
$ws.Cells.Item(130, 26).Value = 0.925
$ws.Cells.Item(130, 27).Value = 0
$ws.Cells.Item(130, 28).Value = 0
$ws.Cells.Item(131, 2).Value = 6815406
$ws.Cells.Item(131, 5).Value = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(131, 6).Value = 'OFK Mladost DG'
$ws.Cells.Item(131, 7).Value = 0
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 9).Value = 'D'
$ws.Cells.Item(131, 10).Value = 2.25
$ws.Cells.Item(131, 11).Value = 3.3
$ws.Cells.Item(131, 12).Value = 2.7
$ws.Cells.Item(131, 13).Value = 2.05
$ws.Cells.Item(131, 14).Value = 3.4
$ws.Cells.Item(131, 15).Value = 3
$ws.Cells.Item(131, 16).Value = -0.25
$ws.Cells.Item(131, 17).Value = 1.825
$ws.Cells.Item(131, 18).Value = 1.975
$ws.Cells.Item(131, 19).Value = 2.5
$ws.Cells.Item(131, 20).Value = 2
$ws.Cells.Item(131, 21).Value = 1.8
$ws.Cells.Item(131, 22).Value = -1
$ws.Cells.Item(131, 23).Value = 2.4
$ws.Cells.Item(131, 25).Value = -0.5
$ws.Cells.Item(131, 26).Value = 0.4875
$ws.Cells.Item(131, 27).Value = -1
$ws.Cells.Item(131, 28).Value = 0.8
$ws.Cells.Item(133, 2).Value = 6815404
$ws.Cells.Item(133, 5).Value = 'Sutjeska Niksic'
$ws.Cells.Item(133, 6).Value = 'FK Rudar Pljevlja'
$ws.Cells.Item(133, 7).Value = 2
$ws.Cells.Item(133, 8).Value = 1
$ws.Cells.Item(133, 9).Value = 'H'
$ws.Cells.Item(133, 10).Value = 1.5
$ws.Cells.Item(133, 11).Value = 3.75
$ws.Cells.Item(133, 12).Value = 5.75
$ws.Cells.Item(133, 13).Value = 1.5
$ws.Cells.Item(133, 14).Value = 3.75
$ws.Cells.Item(133, 15).Value = 5.75
$ws.Cells.Item(133, 16).Value = -1
$ws.Cells.Item(133, 17).Value = 1.85
$ws.Cells.Item(133, 18).Value = 1.95
$ws.Cells.Item(133, 19).Value = 2.25
$ws.Cells.Item(133, 20).Value = 1.9
$ws.Cells.Item(133, 21).Value = 1.9
$ws.Cells.Item(133, 22).Value = 0.5
$ws.Cells.Item(133, 24).Value = -1
$ws.Cells.Item(133, 25).Value = 0
$ws.Cells.Item(133, 26).Value = 0
$ws.Cells.Item(133, 27).Value = 0.8999999999999999
$ws.Cells.Item(133, 28).Value = -1
$ws.Cells.Item(140, 2).Value = 8062092
$ws.Cells.Item(140, 5).Value = 'Sutjeska Niksic'
$ws.Cells.Item(140, 6).Value = 'FK Mornar Bar'
$ws.Cells.Item(140, 7).Value = 2
$ws.Cells.Item(140, 8).Value = 0
$ws.Cells.Item(140, 9).Value = 'H'
$ws.Cells.Item(140, 10).Value = 1.65
$ws.Cells.Item(140, 11).Value = 3.2
$ws.Cells.Item(140, 12).Value = 5
$ws.Cells.Item(140, 13).Value = 1.8
$ws.Cells.Item(140, 14).Value = 3
$ws.Cells.Item(140, 15).Value = 4.5
$ws.Cells.Item(140, 16).Value = -0.5
$ws.Cells.Item(140, 17).Value = 1.825
$ws.Cells.Item(140, 18).Value = 1.975
$ws.Cells.Item(140, 19).Value = 1.75
$ws.Cells.Item(140, 20).Value = 1.775
$ws.Cells.Item(140, 21).Value = 2.025
$ws.Cells.Item(140, 22).Value = 0.8
$ws.Cells.Item(140, 23).Value = -1
$ws.Cells.Item(140, 25).Value = 0.825
$ws.Cells.Item(140, 26).Value = -1
$ws.Cells.Item(140, 27).Value = 0.3875
$ws.Cells.Item(140, 28).Value = -0.5
$ws.Cells.Item(141, 2).Value = 8062093
$ws.Cells.Item(141, 5).Value = 'FK Jezero'
$ws.Cells.Item(141, 6).Value = 'FK Arsenal'
$ws.Cells.Item(141, 7).Value = 4
$ws.Cells.Item(141, 10).Value = 2.1
$ws.Cells.Item(141, 11).Value = 3
$ws.Cells.Item(141, 12).Value = 3.25
$ws.Cells.Item(141, 13).Value = 2.1
$ws.Cells.Item(141, 15).Value = 3.2
$ws.Cells.Item(141, 16).Value = -0.25
$ws.Cells.Item(141, 17).Value = 1.875
$ws.Cells.Item(141, 18).Value = 1.925
$ws.Cells.Item(141, 19).Value = 2.25
$ws.Cells.Item(141, 20).Value = 1.95
$ws.Cells.Item(141, 21).Value = 1.85
$ws.Cells.Item(141, 22).Value = 1.1
$ws.Cells.Item(141, 25).Value = 0.875
$ws.Cells.Item(141, 27).Value = 0.95
$ws.Cells.Item(141, 28).Value = -1
$ws.Cells.Item(142, 2).Value = 8062094
$ws.Cells.Item(142, 5).Value = 'FK Rudar Pljevlja'
$ws.Cells.Item(142, 6).Value = 'FK Jedinstvo Bijelo Polje'
$ws.Cells.Item(142, 7).Value = 1
$ws.Cells.Item(142, 8).Value = 1
$ws.Cells.Item(142, 9).Value = 'D'
$ws.Cells.Item(142, 10).Value = 2.25
$ws.Cells.Item(142, 12).Value = 3
$ws.Cells.Item(142, 13).Value = 2.25
$ws.Cells.Item(142, 14).Value = 3.1
$ws.Cells.Item(142, 15).Value = 2.875
$ws.Cells.Item(142, 17).Value = 2
$ws.Cells.Item(142, 18).Value = 1.8
$ws.Cells.Item(142, 22).Value = -1
$ws.Cells.Item(142, 23).Value = 2.1
$ws.Cells.Item(142, 25).Value = -0.5
$ws.Cells.Item(142, 26).Value = 0.4
$ws.Cells.Item(142, 27).Value = -0.5
$ws.Cells.Item(142, 28).Value = 0.425
